$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The author deleted the row containing "Profanity" (row 20), which shifts
# all subsequent rows up by one.
$ws.Rows.Item(20).Delete() | Out-Null

# Reflect the author's final cell selection recorded in the saved file.
$ws.Range("D39").Select() | Out-Null
